$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date/Week-like text columns auto-convert to numbers/dates when assigned
# directly, so force text entry then strip the resulting number-format
# style so the saved cell carries no explicit style (matching row 2).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-01-02"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "08:58:16"
$ws.Range("C3").Value = "Thursday"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "00"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = 127657
$ws.Range("F3").Value = 143609
$ws.Range("G3").Value = 166391
$ws.Range("H3").Value = 155106
$ws.Range("I3").Value = -1
$ws.Range("J3").Value = 140950
$ws.Range("K3").Value = -1
$ws.Range("L3").Value = -1
$ws.Range("M3").Value = 191198
$ws.Range("N3").Value = 114035
$ws.Range("O3").Value = 45009
$ws.Range("P3").Value = 28129
$ws.Range("Q3").Value = 61883
$ws.Range("R3").Value = -1
$ws.Range("S3").Value = 47754
$ws.Range("T3").Value = -1
